$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("ODI Batting Extra")
$new = $wb.Worksheets.Add()
$new.Name = "ODI Bowling Extra"

$new.Range("A1").Value = "MATCH_CODE"
$new.Range("B1").Value = "MAIDEN_OVERS"
$new.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"
$new.Range("A1:C1").Style = $src.Range("A1:F1").Style
